# Feb 2 2023 22:46 UTC GitHub Actions symbol-list refresh:
# updates the Price (column D) and Volume(1h) (column E) for each coin row.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force text formatting first so Excel keeps the literal strings
# ("327.92", "3.93%", ...) instead of auto-converting them to
# numbers/percentages, matching the source sheet's text cells.
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "327.92"
$ws.Range("E2").Value = "3.93%"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "40.16"
$ws.Range("E3").Value = "6.11%"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.801"
$ws.Range("E4").Value = "11.82%"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08022"
$ws.Range("E5").Value = "0.38%"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "4.574"
$ws.Range("E6").Value = "1.92%"

$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "8.709"
$ws.Range("E7").Value = "1.91%"

$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = "1.938"
$ws.Range("E8").Value = "0.29%"

$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "2.941"
$ws.Range("E9").Value = "-0.74%"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9451"
$ws.Range("E10").Value = "0.06%"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1253"
$ws.Range("E11").Value = "-4.23%"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1965"
$ws.Range("E12").Value = "1.44%"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "8.851"
$ws.Range("E13").Value = "33.57%"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09194"
$ws.Range("E14").Value = "1.18%"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.03589"
$ws.Range("E15").Value = "5.34%"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = "0.09630"
$ws.Range("E16").Value = "1.01%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-6.42%"

$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006246"
$ws.Range("E18").Value = "6.25%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.80%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.34%"

$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1435"
$ws.Range("E21").Value = "10.47%"

$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2418"
$ws.Range("E22").Value = "-0.09%"

$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04399"
$ws.Range("E23").Value = "0.40%"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001260"
$ws.Range("E24").Value = "2.55%"

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004319"
$ws.Range("E25").Value = "1.38%"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001143"
$ws.Range("E26").Value = "-13.96%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.25%"

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02412"
$ws.Range("E39").Value = "0.28%"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05269"
$ws.Range("E40").Value = "2.31%"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007486"
$ws.Range("E41").Value = "-1.65%"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1422"
$ws.Range("E42").Value = "1.60%"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008700"
$ws.Range("E43").Value = "1.27%"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002105"
$ws.Range("E44").Value = "-0.15%"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009481"
$ws.Range("E45").Value = "8.38%"

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006906"
$ws.Range("E46").Value = "6.49%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.61%"

$ws.Range("D48:E48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003158"
$ws.Range("E48").Value = "10.26%"

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001426"
$ws.Range("E49").Value = "-15.43%"

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002109"
$ws.Range("E50").Value = "0.61%"

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002009"
$ws.Range("E51").Value = "0.61%"
